$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "1.00", "238.00")
# are preserved exactly as text instead of being coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "54.364.06"
$ws.Range("E2").Value = "  +1.79%  "

# Row 3
$ws.Range("D3").Value = "2.285.65"
$ws.Range("E3").Value = "  +3.22%  "

# Row 4
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "495.37"
$ws.Range("E5").Value = "  +2.77%  "

# Row 6
$ws.Range("D6").Value = "128.23"
$ws.Range("E6").Value = "  +2.89%  "

# Row 7
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").Value = "0.530"
$ws.Range("E8").Value = "  +2.66%  "

# Row 9
$ws.Range("D9").Value = "2.284.76"
$ws.Range("E9").Value = "  +2.94%  "

# Row 10
$ws.Range("E10").Value = "  +4.48%  "

# Row 11
$ws.Range("E11").Value = "  +2.29%  "

# Row 12
$ws.Range("E12").Value = "  +4.47%  "

# Row 13
$ws.Range("E13").Value = "  -0.23%  "

# Row 14
$ws.Range("D14").Value = "2.689.79"
$ws.Range("E14").Value = "  +2.82%  "

# Row 15
$ws.Range("D15").Value = "21.79"
$ws.Range("E15").Value = "  +4.47%  "

# Row 16
$ws.Range("D16").Value = "54.301.54"
$ws.Range("E16").Value = "  +1.75%  "

# Row 17
$ws.Range("E17").Value = "  +1.88%  "

# Row 18
$ws.Range("D18").Value = "2.297.45"
$ws.Range("E18").Value = "  +3.25%  "

# Row 19
$ws.Range("E19").Value = "  +5.57%  "

# Row 20
$ws.Range("D20").Value = "4.11"
$ws.Range("E20").Value = "  +4.35%  "

# Row 21
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "6.45"
$ws.Range("E21").Value = "  +6.17%  "

# Row 22
$ws.Range("B22").Value = "BitcoinCash"
$ws.Range("C22").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D22").Value = "301.27"
$ws.Range("E22").Value = "  +1.40%  "

# Row 23
$ws.Range("E23").Value = "  +0.11%  "

# Row 24
$ws.Range("E24").Value = "  -2.02%  "

# Row 25
$ws.Range("D25").Value = "62.67"
$ws.Range("E25").Value = "  -0.71%  "

# Row 26
$ws.Range("E26").Value = "  +1.64%  "

# Row 27
$ws.Range("E27").Value = "  +3.20%  "

# Row 28
$ws.Range("D28").Value = "2.385.97"
$ws.Range("E28").Value = "  +2.52%  "

# Row 29
$ws.Range("D29").Value = "0.148"
$ws.Range("E29").Value = "  +4.43%  "

# Row 30
$ws.Range("D30").Value = "7.06"
$ws.Range("E30").Value = "  +2.04%  "

# Row 31
$ws.Range("D31").Value = "169.19"
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("D32").Value = "1.60"
$ws.Range("E32").Value = "  +2.51%  "

# Row 33
$ws.Range("E33").Value = "  +2.79%  "

# Row 34
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  +2.52%  "

# Row 35
$ws.Range("E35").Value = "  +0.09%  "

# Row 36
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.13%  "

# Row 37
$ws.Range("E37").Value = "  +2.85%  "

# Row 38
$ws.Range("D38").Value = "17.68"
$ws.Range("E38").Value = "  +1.92%  "

# Row 39
$ws.Range("D39").Value = "1.19"
$ws.Range("E39").Value = "  +4.76%  "

# Row 40
$ws.Range("D40").Value = "0.880"
$ws.Range("E40").Value = "  +6.59%  "

# Row 41
$ws.Range("E41").Value = "  +4.70%  "

# Row 42
$ws.Range("D42").Value = "35.48"
$ws.Range("E42").Value = "  -0.68%  "

# Row 43
$ws.Range("E43").Value = "  +4.07%  "

# Row 44
$ws.Range("E44").Value = "  +2.86%  "

# Row 45
$ws.Range("E45").Value = "  +3.37%  "

# Row 46
$ws.Range("D46").Value = "127.33"
$ws.Range("E46").Value = "  +4.22%  "

# Row 47
$ws.Range("E47").Value = "  +4.23%  "

# Row 48
$ws.Range("D48").Value = "0.0887"
$ws.Range("E48").Value = "  +1.32%  "

# Row 49
$ws.Range("D49").Value = "0.545"
$ws.Range("E49").Value = "  +2.74%  "

# Row 50
$ws.Range("D50").Value = "238.00"
$ws.Range("E50").Value = "  +4.32%  "

# Row 51
$ws.Range("E51").Value = "  +3.76%  "

# Restore normal style on column D so no stray text-format style lingers
$ws.Range("D2:D51").Style = "Normal"

Write-Output "Applied cryptos list update"